$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 5

# --- Row 4 updates ---
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("X4").Value = 7
$ws.Range("AC4").Value = 6
$ws.Range("AE4").Value = 21
$ws.Range("AF4").Value = 81
$ws.Range("AN4").Value = 3.5
$ws.Range("AO4").Value = 10
$ws.Range("AW4").Value = 6.5

# --- Row 9 updates ---
$ws.Range("U9").Value = 1.83
$ws.Range("V9").Value = 1.83

# --- Row 10 updates ---
$ws.Range("G10").Value = 2.2
$ws.Range("I10").Value = 3.6
$ws.Range("J10").Value = 3
$ws.Range("O10").Value = 1.44
$ws.Range("P10").Value = 2.63
$ws.Range("AE10").Value = 19
$ws.Range("AG10").Value = 8.5
$ws.Range("AH10").Value = 17
$ws.Range("AQ10").Value = 41
$ws.Range("AZ10").Value = 81
$ws.Range("BA10").Value = 126
$ws.Range("BB10").Value = 351

# --- Delete row 11 (Ecuador Liga Pro - Imbabura vs Dep. Cuenca); row 12 (Paraguay) shifts up to 11 ---
$ws.Rows.Item(11).Delete()
